# Updates cryptocurrency price/volume data (and a few re-ranked coin rows)
# to match the latest scrape, per commit "Updated symbol list on Wed Jan 25
# 06:50:07 UTC 2023 with GitHub Actions".
#
# Columns D (Price) and E (Volume(1h)) hold numeric-looking text (e.g.
# "302.71", "-5.04%") that must stay stored as literal text, not be
# reinterpreted as numbers/percentages. We force text by prefixing the
# value with a leading apostrophe (Excel's "store as text" convention,
# which is not itself part of the stored value) and then reset the
# cell's Style back to "Normal" so no stray number-format/quote-prefix
# style lingers on the cell.

function Set-TextValue($Range, $Text) {
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Set-TextValue $ws.Range("D2") '302.71'
Set-TextValue $ws.Range("E2") '-5.04%'
Set-TextValue $ws.Range("D3") '34.95'
Set-TextValue $ws.Range("E3") '-3.28%'
Set-TextValue $ws.Range("D4") '5.061'
Set-TextValue $ws.Range("E4") '-1.39%'
Set-TextValue $ws.Range("D5") '0.07970'
Set-TextValue $ws.Range("E5") '-3.05%'
Set-TextValue $ws.Range("D6") '1.935'
Set-TextValue $ws.Range("E6") '-10.18%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D7") '4.048'
Set-TextValue $ws.Range("E7") '-2.15%'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws.Range("D8") '7.741'
Set-TextValue $ws.Range("E8") '-3.35%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range("D9") '2.955'
Set-TextValue $ws.Range("E9") '5.49%'
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D10") '0.9218'
Set-TextValue $ws.Range("E10") '-0.48%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range("D11") '0.1223'
Set-TextValue $ws.Range("E11") '21.69%'
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D12") '0.1839'
Set-TextValue $ws.Range("E12") '-2.98%'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range("D13") '0.09343'
Set-TextValue $ws.Range("E13") '1.29%'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D14") '0.03524'
Set-TextValue $ws.Range("E14") '-2.28%'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D15") '0.09845'
Set-TextValue $ws.Range("E15") '-0.87%'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D16") '0.001391'
Set-TextValue $ws.Range("E16") '-3.60%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D17") '0.005829'
Set-TextValue $ws.Range("E17") '2.49%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D18") '3.497'
Set-TextValue $ws.Range("E18") '0.95%'
Set-TextValue $ws.Range("D19") '0.3446'
Set-TextValue $ws.Range("E19") '2.16%'
Set-TextValue $ws.Range("D20") '0.1290'
Set-TextValue $ws.Range("E20") '-0.89%'
Set-TextValue $ws.Range("D21") '5.032'
Set-TextValue $ws.Range("E21") '-0.55%'
Set-TextValue $ws.Range("D22") '0.2464'
Set-TextValue $ws.Range("E22") '12.48%'
Set-TextValue $ws.Range("D23") '0.04496'
Set-TextValue $ws.Range("E23") '-2.17%'
Set-TextValue $ws.Range("D25") '0.004854'
Set-TextValue $ws.Range("E25") '2.56%'
Set-TextValue $ws.Range("E26") '-0.08%'
Set-TextValue $ws.Range("E27") '-6.93%'
Set-TextValue $ws.Range("D39") '0.01915'
Set-TextValue $ws.Range("E39") '-4.80%'
Set-TextValue $ws.Range("E40") '-4.85%'
Set-TextValue $ws.Range("D41") '0.007555'
Set-TextValue $ws.Range("E41") '-3.19%'
Set-TextValue $ws.Range("D42") '0.009548'
Set-TextValue $ws.Range("E42") '26.79%'
Set-TextValue $ws.Range("D43") '0.1328'
Set-TextValue $ws.Range("E43") '-5.25%'
Set-TextValue $ws.Range("D44") '0.002109'
Set-TextValue $ws.Range("E44") '0.57%'
Set-TextValue $ws.Range("E45") '-5.83%'
Set-TextValue $ws.Range("D46") '0.00006269'
Set-TextValue $ws.Range("E46") '-3.10%'
Set-TextValue $ws.Range("D47") '0.00000000750'
Set-TextValue $ws.Range("E47") '-0.16%'
Set-TextValue $ws.Range("E49") '-31.43%'
Set-TextValue $ws.Range("D50") '0.00002099'
Set-TextValue $ws.Range("E50") '-0.16%'
Set-TextValue $ws.Range("E51") '-0.16%'

Write-Host "Applied 92 cell updates to match the refreshed crypto symbol list"
